$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 16:48:33"
$ws.Range("E3").Value = "2026-02-20 16:48:35"
$ws.Range("K3").Value = "8.2 MJ/m2"
$ws.Range("O3").Value = "-5.2 °C"
$ws.Range("E4").Value = "2026-02-20 16:48:38"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "52%"
$ws.Range("J4").Value = "1021.8 hPa"
$ws.Range("K4").Value = "7.8 MJ/m2"
$ws.Range("O4").Value = "10.7 °C"
$ws.Range("E5").Value = "2026-02-20 16:48:40"
$ws.Range("K5").Value = "10.4 MJ/m2"
$ws.Range("O5").Value = "-4.6 °C"
$ws.Range("E6").Value = "2026-02-20 16:48:43"
$ws.Range("J6").Value = "1021.8 hPa"
$ws.Range("O6").Value = "9.5 °C"
$ws.Range("E7").Value = "2026-02-20 16:48:45"
$ws.Range("J7").Value = "1021.7 hPa"
$ws.Range("K7").Value = "11.8 MJ/m2"
$ws.Range("O7").Value = "13.3 °C"
$ws.Range("E8").Value = "2026-02-20 16:48:48"
$ws.Range("K8").Value = "10.7 MJ/m2"
$ws.Range("O8").Value = "8.9 °C"
$ws.Range("E9").Value = "2026-02-20 16:48:50"
$ws.Range("K9").Value = "10.8 MJ/m2"
$ws.Range("O9").Value = "13.4 °C"
$ws.Range("E10").Value = "2026-02-20 16:48:53"
$ws.Range("K10").Value = "11.0 MJ/m2"
$ws.Range("O10").Value = "7.9 °C"
$ws.Range("E11").Value = "2026-02-20 16:48:55"
$ws.Range("O11").Value = "9.5 °C"
$ws.Range("E12").Value = "2026-02-20 16:48:58"
$ws.Range("O12").Value = "13.5 °C"
$ws.Range("E13").Value = "2026-02-20 16:49:00"
$ws.Range("K13").Value = "12.5 MJ/m2"
$ws.Range("O13").Value = "6.6 °C"
$ws.Range("E14").Value = "2026-02-20 16:49:03"
$ws.Range("K14").Value = "11.2 MJ/m2"
$ws.Range("O14").Value = "12.3 °C"
$ws.Range("E15").Value = "2026-02-20 16:49:05"
$ws.Range("E16").Value = "2026-02-20 16:49:07"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "48%"
$ws.Range("K16").Value = "11.1 MJ/m2"
$ws.Range("O16").Value = "-3.8 °C"
$ws.Range("E17").Value = "2026-02-20 16:49:10"
$ws.Range("K17").Value = "5.1 MJ/m2"
$ws.Range("O17").Value = "2.6 °C"
$ws.Range("E18").Value = "2026-02-20 16:49:12"
$ws.Range("J18").Value = "1022.1 hPa"
$ws.Range("K18").Value = "9.9 MJ/m2"
$ws.Range("O18").Value = "7.8 °C"
$ws.Range("E19").Value = "2026-02-20 16:49:15"
$ws.Range("K19").Value = "7.1 MJ/m2"
$ws.Range("E20").Value = "2026-02-20 16:49:17"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "55%"
$ws.Range("K20").Value = "14.4 MJ/m2"
$ws.Range("O20").Value = "-3.2 °C"
$ws.Range("E21").Value = "2026-02-20 16:49:20"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "36%"
$ws.Range("K21").Value = "12.6 MJ/m2"
$ws.Range("O21").Value = "9.3 °C"
$ws.Range("E22").Value = "2026-02-20 16:49:22"
$ws.Range("K22").Value = "14.0 MJ/m2"
$ws.Range("E23").Value = "2026-02-20 16:49:25"
$ws.Range("K23").Value = "15.7 MJ/m2"
$ws.Range("E24").Value = "2026-02-20 16:49:27"
$ws.Range("K24").Value = "13.7 MJ/m2"
$ws.Range("O24").Value = "9.3 °C"
$ws.Range("E25").Value = "2026-02-20 16:49:30"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "49%"
$ws.Range("K25").Value = "14.3 MJ/m2"
$ws.Range("O25").Value = "-1.9 °C"
$ws.Range("E26").Value = "2026-02-20 16:49:32"
$ws.Range("J26").Value = "1020.9 hPa"
$ws.Range("K26").Value = "8.9 MJ/m2"
$ws.Range("O26").Value = "5.5 °C"
$ws.Range("E27").Value = "2026-02-20 16:49:35"
$ws.Range("K27").Value = "13.2 MJ/m2"
$ws.Range("O27").Value = "-0.9 °C"
$ws.Range("E28").Value = "2026-02-20 16:49:38"
$ws.Range("J28").Value = "1022.2 hPa"
$ws.Range("K28").Value = "8.0 MJ/m2"
$ws.Range("O28").Value = "7.1 °C"
$ws.Range("E29").Value = "2026-02-20 16:49:40"
$ws.Range("K29").Value = "11.7 MJ/m2"
$ws.Range("O29").Value = "9.2 °C"
$ws.Range("E30").Value = "2026-02-20 16:49:43"
$ws.Range("K30").Value = "10.5 MJ/m2"
$ws.Range("E31").Value = "2026-02-20 16:49:45"
$ws.Range("J31").Value = "1020.6 hPa"
$ws.Range("K31").Value = "13.2 MJ/m2"
$ws.Range("E32").Value = "2026-02-20 16:49:48"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "82%"
$ws.Range("K32").Value = "13.3 MJ/m2"
$ws.Range("O32").Value = "4.3 °C"
$ws.Range("E33").Value = "2026-02-20 16:49:51"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "37%"
$ws.Range("K33").Value = "12.9 MJ/m2"
$ws.Range("O33").Value = "6.3 °C"
$ws.Range("E34").Value = "2026-02-20 16:49:53"
$ws.Range("O34").Value = "0.2 °C"
$ws.Range("E35").Value = "2026-02-20 16:49:56"
$ws.Range("K35").Value = "10.6 MJ/m2"
$ws.Range("O35").Value = "3.7 °C"
$ws.Range("E36").Value = "2026-02-20 16:49:59"
$ws.Range("K36").Value = "12.2 MJ/m2"
$ws.Range("E37").Value = "2026-02-20 16:50:01"
$ws.Range("J37").Value = "1023.6 hPa"
$ws.Range("O37").Value = "5.0 °C"
$ws.Range("E38").Value = "2026-02-20 16:50:04"
$ws.Range("K38").Value = "9.1 MJ/m2"
$ws.Range("O38").Value = "8.9 °C"
$ws.Range("E39").Value = "2026-02-20 16:50:06"
$ws.Range("K39").Value = "14.7 MJ/m2"
$ws.Range("O39").Value = "-3.1 °C"
$ws.Range("E40").Value = "2026-02-20 16:50:09"
$ws.Range("O40").Value = "10.4 °C"
$ws.Range("E41").Value = "2026-02-20 16:50:12"
$ws.Range("J41").Value = "1022.5 hPa"
$ws.Range("K41").Value = "14.1 MJ/m2"
$ws.Range("O41").Value = "13.2 °C"
$ws.Range("E42").Value = "2026-02-20 16:50:14"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "69%"
$ws.Range("O42").Value = "9.6 °C"
$ws.Range("E43").Value = "2026-02-20 16:50:17"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "75%"
$ws.Range("K43").Value = "6.8 MJ/m2"
$ws.Range("O43").Value = "4.7 °C"
$ws.Range("E44").Value = "2026-02-20 16:50:19"
$ws.Range("K44").Value = "9.6 MJ/m2"
$ws.Range("O44").Value = "-5.2 °C"
$ws.Range("E45").Value = "2026-02-20 16:50:22"
$ws.Range("K45").Value = "8.6 MJ/m2"
$ws.Range("E46").Value = "2026-02-20 16:50:25"
$ws.Range("K46").Value = "12.2 MJ/m2"
$ws.Range("O46").Value = "11.9 °C"
